$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row 61 - 10/16/2019: CRM value off, titrator read "Status Not Ok"
# (commit: "10/16/2019 CRM value off - titrator read "Status Not Ok"")

# A61: same date format as the rows above it (m/d/yyyy via style of A60)
$ws.Range("A60").Copy()
$ws.Range("A61").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A61").Value = 43754         # 10/16/2019

$ws.Range("B61").Value = 3167.9356699999998   # Batch value
$ws.Range("C61").Value = 2207.0300000000002   # CRM value
$ws.Range("D61").Formula = "=100*(B61-C61)/C61"   # % off
$ws.Range("E61").Value = 169                       # Batch #
$ws.Range("F61").Value = "junk 5 of 6 ""Not Ok"".  Junk 6 and CRM both inconsistent"

# Scroll the view down to the newly added row and select it, same as the author leaving off there
$excel.ActiveWindow.ScrollRow = 51
$ws.Range("F61").Select()
